# Week 16 stat logging + Week 17 season-sim roster/stat update
# for the Broncos "Players Data" workbook (Rushing + Receiving sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# D.Lock moves ahead of T.Bridgewater in the roster order (Lock is now
# the starter) - swap the two players' name labels on rows 2 and 3,
# carrying each player's own updated counting stats with them.
$rushing.Range("B2").Value = "D.Lock"
$rushing.Range("C2").Value = 0
$rushing.Range("D2").Value = 1
$rushing.Range("E2").Value = 2
$rushing.Range("F2").Value = 1

$rushing.Range("B3").Value = "T.Bridgewater"
$rushing.Range("C3").Value = 3
$rushing.Range("D3").Value = 7
$rushing.Range("E3").Value = 13
$rushing.Range("F3").Value = 8

# M.Gordon
$rushing.Range("C4").Value = 98
$rushing.Range("D4").Value = 69

# J.Williams
$rushing.Range("C5").Value = 96
$rushing.Range("D5").Value = 60
$rushing.Range("F5").Value = 22

# ---------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# M.Gordon
$receiving.Range("C2").Value = 31
$receiving.Range("D2").Value = 22

# J.Williams
$receiving.Range("C3").Value = 48
$receiving.Range("D3").Value = 38

# C.Sutton
$receiving.Range("C5").Value = 59
$receiving.Range("D5").Value = 46
$receiving.Range("E5").Value = 31

# D.Hamilton
$receiving.Range("C6").Value = 75
$receiving.Range("D6").Value = 59
$receiving.Range("E6").Value = 23
$receiving.Range("F6").Value = 17

# T.Patrick
$receiving.Range("C7").Value = 57
$receiving.Range("D7").Value = 39
$receiving.Range("E7").Value = 18

# K.Hinton
$receiving.Range("C10").Value = 74
$receiving.Range("D10").Value = 59
$receiving.Range("E10").Value = 13
$receiving.Range("F10").Value = 8

# N.Fant
$receiving.Range("E11").Value = 7
